# TC37_Canine_Filter_Breed-ParsonRussTerr.xlsx
# "updated first 25 tc in icdc breed+diagnosis"
#
# The workbook's first sheet is internally named "startup" (it is the sheet
# that carries xl/worksheets/sheet1.xml / rId1) and holds the TabName /
# query / StatQuery / dbExcel / WebExcel lookup table consumed by the rest
# of the automation. "Sheet1" (xl/worksheets/sheet2.xml) is just an empty
# helper tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- B2: CasesTab "query" column -> append the Cohort coalesce line ------
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Parson Russell Terrier'] 
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@
$ws.Range("B2").Value = $casesQuery

# --- B4: FilesTab "query" column -> drop the trailing Study Code line ----
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Parson Russell Terrier'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis  
'@
$ws.Range("B4").Value = $filesQuery

# --- C2/C3/C4: "StatQuery"/dbExcel column -> replace the old
# all_studies/all_breeds rollup query with the new Programs/Studies/Cases/
# Samples/Case Files/Study Files rollup query (same text in all 3 rows). --
$dbExcelQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Parson Russell Terrier'] 
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$ws.Range("C2").Value = $dbExcelQuery
$ws.Range("C3").Value = $dbExcelQuery
$ws.Range("C4").Value = $dbExcelQuery

# --- Row heights / column widths --------------------------------------
# The workbook was re-saved from a newer Excel build, which re-wrapped the
# long StatQuery/query text at slightly different metrics and bumped the
# sheet's row heights / column widths a little. Reproduce the row heights
# (exact) and column widths (nearest reachable) from the target file.
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

$ws.Columns.Item(1).ColumnWidth = 10.0221354166667
$ws.Columns.Item(2).ColumnWidth = 91.5924479166667
$ws.Columns.Item(3).ColumnWidth = 85.4518229166667
$ws.Columns.Item(4).ColumnWidth = 69.4518229166667
$ws.Columns.Item(5).ColumnWidth = 39.8776041666667
